{"js": "// Update the date line and the 25 multiplication problems in the table,\n// matching the author's commit that regenerated this worksheet's numbers.\nconst replacements = [\n  [\"2024-08-31 Saturday\", \"2024-09-01 Sunday\"],\n  [\"612\u00d73=1836\", \"120\u00d73=360\"],\n  [\"528\u00d77=3696\", \"386\u00d73=1158\"],\n  [\"653\u00d72=1306\", \"382\u00d78=3056\"],\n  [\"469\u00d73=1407\", \"838\u00d74=3352\"],\n  [\"143\u00d75=715\", \"476\u00d78=3808\"],\n  [\"503\u00d74=2012\", \"732\u00d76=4392\"],\n  [\"562\u00d79=5058\", \"530\u00d78=4240\"],\n  [\"376\u00d79=3384\", \"888\u00d79=7992\"],\n  [\"774\u00d74=3096\", \"209\u00d72=418\"],\n  [\"982\u00d75=4910\", \"563\u00d75=2815\"],\n  [\"483\u00d75=2415\", \"856\u00d73=2568\"],\n  [\"110\u00d73=330\", \"161\u00d79=1449\"],\n  [\"966\u00d76=5796\", \"415\u00d76=2490\"],\n  [\"113\u00d79=1017\", \"225\u00d74=900\"],\n  [\"660\u00d78=5280\", \"843\u00d72=1686\"],\n  [\"217\u00d76=1302\", \"355\u00d72=710\"],\n  [\"697\u00d72=1394\", \"440\u00d75=2200\"],\n  [\"848\u00d77=5936\", \"643\u00d73=1929\"],\n  [\"538\u00d72=1076\", \"434\u00d78=3472\"],\n  [\"848\u00d73=2544\", \"497\u00d72=994\"],\n  [\"513\u00d74=2052\", \"937\u00d72=1874\"],\n  [\"700\u00d77=4900\", \"240\u00d75=1200\"],\n  [\"945\u00d74=3780\", \"340\u00d79=3060\"],\n  [\"499\u00d78=3992\", \"901\u00d74=3604\"],\n  [\"905\u00d79=8145\", \"832\u00d76=4992\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication problems in the table,\n# matching the author's commit that regenerated this worksheet's numbers.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-31 Saturday\", \"2024-09-01 Sunday\"),\n    @(\"612\u00d73=1836\", \"120\u00d73=360\"),\n    @(\"528\u00d77=3696\", \"386\u00d73=1158\"),\n    @(\"653\u00d72=1306\", \"382\u00d78=3056\"),\n    @(\"469\u00d73=1407\", \"838\u00d74=3352\"),\n    @(\"143\u00d75=715\", \"476\u00d78=3808\"),\n    @(\"503\u00d74=2012\", \"732\u00d76=4392\"),\n    @(\"562\u00d79=5058\", \"530\u00d78=4240\"),\n    @(\"376\u00d79=3384\", \"888\u00d79=7992\"),\n    @(\"774\u00d74=3096\", \"209\u00d72=418\"),\n    @(\"982\u00d75=4910\", \"563\u00d75=2815\"),\n    @(\"483\u00d75=2415\", \"856\u00d73=2568\"),\n    @(\"110\u00d73=330\", \"161\u00d79=1449\"),\n    @(\"966\u00d76=5796\", \"415\u00d76=2490\"),\n    @(\"113\u00d79=1017\", \"225\u00d74=900\"),\n    @(\"660\u00d78=5280\", \"843\u00d72=1686\"),\n    @(\"217\u00d76=1302\", \"355\u00d72=710\"),\n    @(\"697\u00d72=1394\", \"440\u00d75=2200\"),\n    @(\"848\u00d77=5936\", \"643\u00d73=1929\"),\n    @(\"538\u00d72=1076\", \"434\u00d78=3472\"),\n    @(\"848\u00d73=2544\", \"497\u00d72=994\"),\n    @(\"513\u00d74=2052\", \"937\u00d72=1874\"),\n    @(\"700\u00d77=4900\", \"240\u00d75=1200\"),\n    @(\"945\u00d74=3780\", \"340\u00d79=3060\"),\n    @(\"499\u00d78=3992\", \"901\u00d74=3604\"),\n    @(\"905\u00d79=8145\", \"832\u00d76=4992\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
